$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10 (pushes existing rows 10-36 down to 11-37)
$ws.Rows("10:10").Insert()

# Populate the newly inserted row 10 with the new vulnerability entry
$ws.Range("A10").Value = "Node"
$ws.Range("B10").Value = "a,a+"
$ws.Range("C10").Value = "PTV-NET-IDENT-ACTIVE-MLDNVERDEV2"
$ws.Range("D10").Value = "Device only responds to illegitimate MLDv1 queries even though MLDv2 queries are sent, possibly downgraded"

# Append two new rows (38, 39) with new vulnerability entries
$ws.Range("A38").Value = "Network"
$ws.Range("B38").Value = "p,a,a+"
$ws.Range("A39").Value = "Node"
$ws.Range("B39").Value = "p,a,a+"

$ws.Range("C38").Value = "PTV-NET-MITM-ICMP6REDIR"
$ws.Range("C39").Value = "PTV-NET-MITM-ICMP6REDIRDEV"

$ws.Range("D38").Value = "Network does not block ICMPv6 Redirect messages"
$ws.Range("D39").Value = "Device communication can be redirected using ICMPv6 Redirect"

# Update the view: move the selection (also resets the scrolled top-left cell)
$ws.Range("K36").Select()
